# Apply the edits described by the commit diff:
#  1. Rename the worksheet "Monster" -> "Monster.xlsx"
#  2. Shrink every data row (1-42) from the 18.75pt default down to 14.25pt
#  3. Swap the theme's accent1 / accent5 colors (4472C4 <-> 5B9BD5)
#  4. Switch the page setup to Portrait orientation

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Sheet name
$ws.Name = "Monster.xlsx"

# 2) Row heights for the used range (rows 1-42)
$ws.Rows("1:42").RowHeight = 14.25

# 3) Theme accent color swap (accent1 <-> accent5)
$accent1 = $wb.Theme.ThemeColorScheme.Colors(5)
$accent5 = $wb.Theme.ThemeColorScheme.Colors(9)
$accent1Rgb = $accent1.RGB
$accent5Rgb = $accent5.RGB
$accent1.RGB = $accent5Rgb
$accent5.RGB = $accent1Rgb

# 4) Page orientation -> portrait
$ws.PageSetup.Orientation = $xlPortrait
